# Atualização de bases das ligas, do dia: 27-04-2024 às 09:20
# Swap the data (columns B through AB) between row 83 and row 84,
# leaving column A (the running id/index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng1 = $ws.Range("B83:AB83")
$rng2 = $ws.Range("B84:AB84")

$vals1 = $rng1.Value()
$vals2 = $rng2.Value()

$rng1.Value = $vals2
$rng2.Value = $vals1
